$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1366.8889
$ws.Range("J17").Value = 1366.8889
$ws.Range("L17").Value = 4100.6667
$ws.Range("N17").Value = -4436.6667

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 480
$ws.Range("I41").Value = 480
$ws.Range("K41").Value = 480
$ws.Range("M41").Value = -40

# ALC row 109
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 63561
$ws.Range("J109").Value = 63561
$ws.Range("L109").Value = 63561
$ws.Range("N109").Value = -66335

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 9728.833000000001
$ws.Range("I111").Value = 10785.8
$ws.Range("J111").Value = 4444
$ws.Range("K111").Value = 32357.4
$ws.Range("L111").Value = 13332
$ws.Range("M111").Value = -29290.4
$ws.Range("N111").Value = -19466

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2573.4736
$ws.Range("I116").Value = 2338.4614
$ws.Range("J116").Value = 3082.6667
$ws.Range("K116").Value = 2338.4614
$ws.Range("L116").Value = 3082.6667
$ws.Range("M116").Value = 1103.5386
$ws.Range("N116").Value = -9966.6667

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4196.1787
$ws.Range("J137").Value = 5087.2354
$ws.Range("L137").Value = 15261.7062
$ws.Range("N137").Value = -20361.7062

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1412.4546
$ws.Range("I2").Value = 1628.6923
$ws.Range("J2").Value = 1100.1111
$ws.Range("K2").Value = 1628.6923
$ws.Range("L2").Value = 1100.1111
$ws.Range("M2").Value = -1515.6923
$ws.Range("N2").Value = -1326.1111

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2675.0952
$ws.Range("I74").Value = 2760.1875
$ws.Range("J74").Value = 2402.8
$ws.Range("K74").Value = 2760.1875
$ws.Range("L74").Value = 2402.8
$ws.Range("M74").Value = -1886.1875
$ws.Range("N74").Value = -4150.8

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2675.0952
$ws.Range("I77").Value = 2760.1875
$ws.Range("J77").Value = 2402.8
$ws.Range("K77").Value = 13800.9375
$ws.Range("L77").Value = 12014
$ws.Range("M77").Value = -9432.9375
$ws.Range("N77").Value = -20750

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1412.4546
$ws.Range("I116").Value = 1628.6923
$ws.Range("J116").Value = 1100.1111
$ws.Range("K116").Value = 1628.6923
$ws.Range("L116").Value = 1100.1111
$ws.Range("M116").Value = 665.3077000000001
$ws.Range("N116").Value = -5688.1111

# ARM row 118
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 38600
$ws.Range("J118").Value = 38600
$ws.Range("L118").Value = 38600
$ws.Range("N118").Value = -41914

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1249.7812
$ws.Range("I122").Value = 1073.9656
$ws.Range("J122").Value = 2949.3333
$ws.Range("K122").Value = 3221.8968
$ws.Range("L122").Value = 8847.999899999999
$ws.Range("M122").Value = -771.8968
$ws.Range("N122").Value = -13747.9999

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 51510.5
$ws.Range("J139").Value = 51510.5
$ws.Range("L139").Value = 51510.5
$ws.Range("N139").Value = -61790.5

# ARM row 140
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 35357.5
$ws.Range("J140").Value = 35357.5
$ws.Range("L140").Value = 35357.5
$ws.Range("N140").Value = -45717.5

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1412.4546
$ws.Range("I3").Value = 1628.6923
$ws.Range("J3").Value = 1100.1111
$ws.Range("K3").Value = 1628.6923
$ws.Range("L3").Value = 1100.1111
$ws.Range("M3").Value = -1514.6923
$ws.Range("N3").Value = -1328.1111

# BSM row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 17686.666
$ws.Range("J81").Value = 17686.666
$ws.Range("L81").Value = 17686.666
$ws.Range("N81").Value = -19808.666

# BSM row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 17686.666
$ws.Range("J84").Value = 17686.666
$ws.Range("L84").Value = 53059.99800000001
$ws.Range("N84").Value = -63667.99800000001

# BSM row 108
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 39684
$ws.Range("J108").Value = 39684
$ws.Range("L108").Value = 39684
$ws.Range("N108").Value = -47364

# BSM row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 73300
$ws.Range("J138").Value = 73300
$ws.Range("L138").Value = 73300
$ws.Range("N138").Value = -83580

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7041.48
$ws.Range("I31").Value = 8541.666999999999
$ws.Range("K31").Value = 8541.666999999999
$ws.Range("M31").Value = -8246.666999999999

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7041.48
$ws.Range("I34").Value = 8541.666999999999
$ws.Range("K34").Value = 8541.666999999999
$ws.Range("M34").Value = -8339.666999999999

# CRP row 44
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 24442.5
$ws.Range("I44").Value = 9999
$ws.Range("J44").Value = 29257
$ws.Range("K44").Value = 9999
$ws.Range("L44").Value = 29257
$ws.Range("M44").Value = -9557
$ws.Range("N44").Value = -30141

# CRP row 139
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 819.8
$ws.Range("I7").Value = 799.75
$ws.Range("K7").Value = 2399.25
$ws.Range("M7").Value = -2287.25

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 159.95
$ws.Range("I23").Value = 46.666668
$ws.Range("J23").Value = 179.94118
$ws.Range("K23").Value = 140.000004
$ws.Range("L23").Value = 539.82354
$ws.Range("M23").Value = 94.99999600000001
$ws.Range("N23").Value = -1009.82354

# CUL row 29
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 344.33334
$ws.Range("J29").Value = 344.33334
$ws.Range("L29").Value = 1033.00002
$ws.Range("N29").Value = -1587.00002

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 61.066666
$ws.Range("J38").Value = 113.333336
$ws.Range("L38").Value = 340.000008
$ws.Range("N38").Value = -1034.000008

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 378.5
$ws.Range("I98").Value = 351.73334
$ws.Range("J98").Value = 780
$ws.Range("K98").Value = 1055.20002
$ws.Range("L98").Value = 2340
$ws.Range("M98").Value = 442.79998
$ws.Range("N98").Value = -5336

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 665.14
$ws.Range("I113").Value = 678.5897
$ws.Range("J113").Value = 617.4545000000001
$ws.Range("K113").Value = 2035.7691
$ws.Range("L113").Value = 1852.3635
$ws.Range("M113").Value = 134.2309
$ws.Range("N113").Value = -6192.3635

# GSM row 47
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 17515.5
$ws.Range("J47").Value = 17515.5
$ws.Range("L47").Value = 17515.5
$ws.Range("N47").Value = -18651.5

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2400.6897
$ws.Range("I132").Value = 2224.6538
$ws.Range("J132").Value = 3926.3333
$ws.Range("K132").Value = 6673.9614
$ws.Range("L132").Value = 11778.9999
$ws.Range("M132").Value = -4143.9614
$ws.Range("N132").Value = -16838.9999

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 13948.8
$ws.Range("I61").Value = 17614
$ws.Range("J61").Value = 5396.6665
$ws.Range("K61").Value = 17614
$ws.Range("L61").Value = 5396.6665
$ws.Range("M61").Value = -17412
$ws.Range("N61").Value = -5800.6665

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4171.263
$ws.Range("I100").Value = 2745
$ws.Range("J100").Value = 5756
$ws.Range("K100").Value = 2745
$ws.Range("L100").Value = 5756
$ws.Range("M100").Value = -2204
$ws.Range("N100").Value = -6838

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 13948.8
$ws.Range("I113").Value = 17614
$ws.Range("J113").Value = 5396.6665
$ws.Range("K113").Value = 17614
$ws.Range("L113").Value = 5396.6665
$ws.Range("M113").Value = -15444
$ws.Range("N113").Value = -9736.666499999999

# LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 56943.668
$ws.Range("J140").Value = 56943.668
$ws.Range("L140").Value = 56943.668
$ws.Range("N140").Value = -67303.66800000001

# WVR row 111
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 49543.668
$ws.Range("J111").Value = 49543.668
$ws.Range("L111").Value = 49543.668
$ws.Range("N111").Value = -57723.668

# WVR row 140
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 58982.25
$ws.Range("J140").Value = 58982.25
$ws.Range("L140").Value = 58982.25
$ws.Range("N140").Value = -69342.25

Write-Output "Applied all Pandaemonium Profits updates"